$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 6565.833333333333, 7112, 6140, 0.1789425134658814),
    @(1, 6427.966666666666, 6829, 5647, 0.178994623819987),
    @(2, 6250.4, 6755, 5387, 0.1807592789332072),
    @(3, 6833.4, 7438, 6475, 0.1782646258672078),
    @(4, 6334.6, 6855, 5532, 0.1855825662612915),
    @(5, 6318.966666666666, 6790, 5859, 0.1798436959584554),
    @(6, 6678.166666666667, 7159, 5812, 0.1846895535786947),
    @(7, 6623.933333333333, 7157, 5969, 0.1812933127085368),
    @(8, 6418.9, 6988, 5754, 0.1814293146133423),
    @(9, 6653.5, 7311, 5966, 0.1801905870437622)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
